# update code tao report luong tai report co so
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new worksheet "Đơn 1 bác sĩ" right after "Đơn sale chính"
#    (i.e. right before "Lương").
# ------------------------------------------------------------------
$donSaleChinh = $wb.Worksheets.Item(1)
$donBacSi = $wb.Worksheets.Add($null, $donSaleChinh)
$donBacSi.Name = "Đơn 1 bác sĩ"

# Header row
$donBacSi.Cells.Item(1, 1).Value = "Tiền tố"
$donBacSi.Cells.Item(1, 2).Value = "Mã dịch vụ"
$donBacSi.Cells.Item(1, 3).Value = "Ngày thực hiện"
$donBacSi.Cells.Item(1, 4).Value = "Cơ sở"
$donBacSi.Cells.Item(1, 5).Value = "Khách hàng"
$donBacSi.Cells.Item(1, 6).Value = "Nguồn khách"
$donBacSi.Cells.Item(1, 7).Value = "Tên dịch vụ"
$donBacSi.Cells.Item(1, 8).Value = "Đơn giá gốc"
$donBacSi.Cells.Item(1, 9).Value = "Sale phụ"
$donBacSi.Cells.Item(1, 10).Value = "Upsale"
$donBacSi.Cells.Item(1, 11).Value = "Đơn giá"
$donBacSi.Cells.Item(1, 12).Value = "Đã thanh toán"
$donBacSi.Cells.Item(1, 13).Value = "Tỉ lệ chiết khấu bác sĩ 1"
$donBacSi.Cells.Item(1, 14).Value = "Chiết khấu bác sĩ 1"

# Row 2
$donBacSi.Cells.Item(2, 1).Value = "HD-LUXURY"
$donBacSi.Cells.Item(2, 2).Value = 614
$donBacSi.Cells.Item(2, 3).NumberFormat = "@"
$donBacSi.Cells.Item(2, 3).Value = "08-01-2024"
$donBacSi.Cells.Item(2, 4).Value = "CẦN THƠ"
$donBacSi.Cells.Item(2, 5).Value = "Trần Nguyễn Yến Linh"
$donBacSi.Cells.Item(2, 6).Value = "Khách cũ"
$donBacSi.Cells.Item(2, 7).Value = "Cắt mí"
$donBacSi.Cells.Item(2, 8).Value = 0
$donBacSi.Cells.Item(2, 9).Value = "Đỗ Thị Huyền Trân"
$donBacSi.Cells.Item(2, 10).Value = 6000000
$donBacSi.Cells.Item(2, 11).Value = 6000000
$donBacSi.Cells.Item(2, 12).Value = 6000000
$donBacSi.Cells.Item(2, 13).Value = 0.08
$donBacSi.Cells.Item(2, 14).Value = 480000

# Row 3
$donBacSi.Cells.Item(3, 1).Value = "HD-LUXURY"
$donBacSi.Cells.Item(3, 2).Value = 615
$donBacSi.Cells.Item(3, 3).NumberFormat = "@"
$donBacSi.Cells.Item(3, 3).Value = "08-01-2024"
$donBacSi.Cells.Item(3, 4).Value = "CẦN THƠ"
$donBacSi.Cells.Item(3, 5).Value = "Nguyễn Thị Mỹ Duyên"
$donBacSi.Cells.Item(3, 6).Value = "Khách cũ"
$donBacSi.Cells.Item(3, 7).Value = "Tiêm Filler"
$donBacSi.Cells.Item(3, 8).Value = 2100000
$donBacSi.Cells.Item(3, 11).Value = 2100000
$donBacSi.Cells.Item(3, 12).Value = 2100000
$donBacSi.Cells.Item(3, 13).Value = 0.08
$donBacSi.Cells.Item(3, 14).Value = 168000

# Row 4
$donBacSi.Cells.Item(4, 1).Value = "HD-LUXURY"
$donBacSi.Cells.Item(4, 2).Value = 616
$donBacSi.Cells.Item(4, 3).NumberFormat = "@"
$donBacSi.Cells.Item(4, 3).Value = "08-02-2024"
$donBacSi.Cells.Item(4, 4).Value = "LONG XUYÊN"
$donBacSi.Cells.Item(4, 5).Value = "Chị duyên"
$donBacSi.Cells.Item(4, 6).Value = "Khách cũ giới thiệu"
$donBacSi.Cells.Item(4, 7).Value = "Cắt mí"
$donBacSi.Cells.Item(4, 8).Value = 8000000
$donBacSi.Cells.Item(4, 11).Value = 8000000
$donBacSi.Cells.Item(4, 12).Value = 8000000
$donBacSi.Cells.Item(4, 13).Value = 0.1
$donBacSi.Cells.Item(4, 14).Value = 800000

# Row 5
$donBacSi.Cells.Item(5, 1).Value = "HD-LUXURY"
$donBacSi.Cells.Item(5, 2).Value = 617
$donBacSi.Cells.Item(5, 3).NumberFormat = "@"
$donBacSi.Cells.Item(5, 3).Value = "08-02-2024"
$donBacSi.Cells.Item(5, 4).Value = "LONG XUYÊN"
$donBacSi.Cells.Item(5, 5).Value = "Cô tú"
$donBacSi.Cells.Item(5, 6).Value = "Khách cũ"
$donBacSi.Cells.Item(5, 7).Value = "Nâng cung chân mày"
$donBacSi.Cells.Item(5, 8).Value = 4000000
$donBacSi.Cells.Item(5, 11).Value = 4000000
$donBacSi.Cells.Item(5, 12).Value = 3000000
$donBacSi.Cells.Item(5, 13).Value = 0.1
$donBacSi.Cells.Item(5, 14).Value = 300000

# Row 6 - totals
$donBacSi.Cells.Item(6, 1).Value = "Tổng"
$donBacSi.Cells.Item(6, 2).Value = 4
$donBacSi.Cells.Item(6, 8).Value = 14100000
$donBacSi.Cells.Item(6, 10).Value = 6000000
$donBacSi.Cells.Item(6, 11).Value = 20100000
$donBacSi.Cells.Item(6, 12).Value = 19100000
$donBacSi.Cells.Item(6, 13).Value = 0
$donBacSi.Cells.Item(6, 14).Value = 1748000

# ------------------------------------------------------------------
# 2) Update the "Lương" worksheet with the new computed figures.
# ------------------------------------------------------------------
$luong = $wb.Worksheets.Item("Lương")

$luong.Cells.Item(2, 2).Value = 2
$luong.Cells.Item(3, 2).Value = 70000
$luong.Cells.Item(4, 2).Value = 571428.5714285715
$luong.Cells.Item(7, 2).Value = 648000
$luong.Cells.Item(12, 2).Value = 1
$luong.Cells.Item(13, 2).Value = 285714.2857142857
$luong.Cells.Item(17, 2).Value = 1100000
$luong.Cells.Item(32, 2).Value = 1289428.571428571
$luong.Cells.Item(33, 2).Value = 1385714.285714286
$luong.Cells.Item(35, 1).Value = "Tổng lương tại HỆ THỐNG"
$luong.Cells.Item(35, 2).Value = 2675142.857142857
